# Generate Report for Handoff
# The file d395d5b3-774c-41e0-b708-6462d2c4b6c6.md is now "Ready for handoff"
# for both the zh-cn and de-de locales. Update the Overview sheet and each
# locale sheet accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
# Row 3 corresponds to d395d5b3-774c-41e0-b708-6462d2c4b6c6.md
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 14:14:54"

# ---- zh-cn sheet ----
# Row 3 corresponds to d395d5b3-774c-41e0-b708-6462d2c4b6c6.md
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-28 14:14:50"

# ---- de-de sheet ----
# Row 3 corresponds to d395d5b3-774c-41e0-b708-6462d2c4b6c6.md
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-28 14:14:54"

# ---- Column widths (auto-adjusted by Excel after the longer status text) ----
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
